# Weekly data refresh: insert a new record row for Alcachofa (Vega Modelo de
# Temuco) above the current row 275, pushing the existing historical rows
# (275-312) down by one (becoming 276-313).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 275 - shifts old rows 275..312 down to 276..313.
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the latest week's data.
$ws.Cells.Item(275, 1).Value = 10
$ws.Cells.Item(275, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(275, 3).Value = "La Araucanía"
$ws.Cells.Item(275, 4).Value = 45124
$ws.Cells.Item(275, 5).Value = 9
$ws.Cells.Item(275, 6).Value = 100112013
$ws.Cells.Item(275, 7).Value = "Alcachofa"
$ws.Cells.Item(275, 8).Value = "Madrigal"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 500
$ws.Cells.Item(275, 11).Value = 12000
$ws.Cells.Item(275, 12).Value = 15000
$ws.Cells.Item(275, 13).Value = 13800
$ws.Cells.Item(275, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(275, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(275, 16).Value = 345
$ws.Cells.Item(275, 17).Value = 40
$ws.Cells.Item(275, 18).Value = "Hortaliza"
